$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab (matches commit renaming "7 değişkenli" -> "veriset_7")
$ws.Name = "veriset_7"

# Rename the header cells: replace spaces with underscores in the multi-word headers
$ws.Range("C1").Value = "İş_Deneyimi"
$ws.Range("D1").Value = "Harcama_Miktarı"
$ws.Range("E1").Value = "Çalışma_Saati"
$ws.Range("F1").Value = "Mutluluk_Skoru"
$ws.Range("G1").Value = "Sağlık_Durumu"

# Update the view: scroll position (topLeftCell -> A10) and active selection (K26)
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("K26").Select()
